$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 9; this shifts rows 9..86 down to 10..87
# and carries the date-format style from column D down with it.
$ws.Rows.Item(9).Insert()

# Populate the new row 9 with the new weekly price record.
$ws.Cells.Item(9, 1).Value2 = 6
$ws.Cells.Item(9, 2).Value2 = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(9, 3).Value2 = "Metropolitana"
$ws.Cells.Item(9, 4).Value2 = 45111
$ws.Cells.Item(9, 5).Value2 = 13
$ws.Cells.Item(9, 6).Value2 = 100112035
$ws.Cells.Item(9, 7).Value2 = "Bruselas (repollito)"
$ws.Cells.Item(9, 8).Value2 = "Sin especificar"
$ws.Cells.Item(9, 9).Value2 = "Primera"
$ws.Cells.Item(9, 10).Value2 = 380
$ws.Cells.Item(9, 11).Value2 = 19000
$ws.Cells.Item(9, 12).Value2 = 20000
$ws.Cells.Item(9, 13).Value2 = 19395
$ws.Cells.Item(9, 14).Value2 = "$/malla 15 kilos"
$ws.Cells.Item(9, 15).Value2 = "Provincia de Quillota"
$ws.Cells.Item(9, 16).Value2 = 1293
$ws.Cells.Item(9, 17).Value2 = 15
$ws.Cells.Item(9, 18).Value2 = "Hortaliza"
